# Super Lig 2023-2024: fix two mis-ordered fixture pairs (home/away data swapped)
# and append 15 newly scraped fixtures (rows 138-152).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: swap the full F:V block between mis-ordered row pairs ---
$swapPairs = @(
    ,(8, 9)
    ,(18, 19)
    ,(29, 30)
    ,(38, 39)
    ,(40, 41)
    ,(49, 50)
    ,(75, 76)
    ,(79, 80)
    ,(127, 128)
)

foreach ($pair in $swapPairs) {
    $rowA = $pair[0]
    $rowB = $pair[1]
    $rangeA = $ws.Range("F$rowA" + ":V$rowA")
    $rangeB = $ws.Range("F$rowB" + ":V$rowB")
    $valsA = $rangeA.Value2
    $valsB = $rangeB.Value2
    $rangeA.Value2 = $valsB
    $rangeB.Value2 = $valsA
}

# --- Step 2: append 15 new fixtures as rows 138-152 (Indice 137-151) ---
# Copy formatting (bold/border/center style on col A, datetime format on col E) from the last existing row
$ws.Range("A137:V137").Copy() | Out-Null
$ws.Range("A138:V152").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

$newRows = @(
    ,(137, "turkey", "super-lig", "2023-2024", 45263.47916666666, "Kasimpasa", 4, "Gaziantep", 2, 2.1, "25/11/2023 14:12", 2.36, "03/12/2023 11:26", 3.69, "25/11/2023 14:12", 3.5, "03/12/2023 11:21", 3.42, "25/11/2023 14:12", 3.15, "03/12/2023 11:26", "https://www.betexplorer.com/football/turkey/super-lig/kasimpasa-gaziantep/xOf41wuU/")
    ,(138, "turkey", "super-lig", "2023-2024", 45263.58333333334, "Trabzonspor", 0, "Kayserispor", 1, 1.72, "27/11/2023 18:12", 1.61, "03/12/2023 13:54", 4.09, "27/11/2023 18:12", 4.29, "03/12/2023 13:54", 4.78, "27/11/2023 18:12", 5.69, "03/12/2023 13:54", "https://www.betexplorer.com/football/turkey/super-lig/trabzonspor-kayserispor/l2OMJKfb/")
    ,(139, "turkey", "super-lig", "2023-2024", 45263.70833333334, "Ankaragucu", 1, "Besiktas", 1, 3.48, "26/11/2023 14:13", 3.39, "03/12/2023 16:55", 3.67, "26/11/2023 14:13", 3.65, "03/12/2023 16:55", 2.12, "26/11/2023 14:13", 2.18, "03/12/2023 16:53", "https://www.betexplorer.com/football/turkey/super-lig/ankaragucu-besiktas/84vOHbPA/")
    ,(140, "turkey", "super-lig", "2023-2024", 45264.75, "Fenerbahce", 4, "Sivasspor", 1, 1.24, "27/11/2023 18:12", 1.19, "04/12/2023 17:53", 6.96, "27/11/2023 18:12", 7.74, "04/12/2023 17:59", 11.91, "27/11/2023 18:12", 15.12, "04/12/2023 17:59", "https://www.betexplorer.com/football/turkey/super-lig/fenerbahce-sivasspor/ShNQIvA4/")
    ,(141, "turkey", "super-lig", "2023-2024", 45268.75, "Galatasaray", 3, "Adana Demirspor", 1, 1.38, "02/12/2023 18:13", 1.27, "08/12/2023 17:54", 5.71, "02/12/2023 18:13", 6.71, "08/12/2023 17:34", 6.84, "02/12/2023 18:13", 10.1, "08/12/2023 17:34", "https://www.betexplorer.com/football/turkey/super-lig/galatasaray-adanademirspor/02J2AGHp/")
    ,(142, "turkey", "super-lig", "2023-2024", 45269.47916666666, "Antalyaspor", 2, "Karagumruk", 1, 1.79, "02/12/2023 18:13", 1.93, "09/12/2023 11:29", 3.75, "02/12/2023 18:13", 3.62, "09/12/2023 11:29", 4.61, "02/12/2023 18:13", 4.23, "09/12/2023 11:29", "https://www.betexplorer.com/football/turkey/super-lig/antalyaspor-f-karagumruk/lOacWcv4/")
    ,(143, "turkey", "super-lig", "2023-2024", 45269.58333333334, "Samsunspor", 4, "Kasimpasa", 2, 2.03, "03/12/2023 11:42", 1.84, "09/12/2023 13:57", 3.86, "03/12/2023 11:42", 4.13, "09/12/2023 13:57", 3.45, "03/12/2023 11:42", 4.08, "09/12/2023 13:57", "https://www.betexplorer.com/football/turkey/super-lig/samsunspor-kasimpasa/vgI69zXj/")
    ,(144, "turkey", "super-lig", "2023-2024", 45269.70833333334, "Besiktas", 1, "Fenerbahce", 3, 3.53, "03/12/2023 20:12", 4.33, "09/12/2023 16:59", 3.56, "03/12/2023 20:12", 3.92, "09/12/2023 16:59", 2.1, "03/12/2023 20:12", 1.84, "09/12/2023 16:59", "https://www.betexplorer.com/football/turkey/super-lig/besiktas-fenerbahce/n9I0VHgA/")
    ,(145, "turkey", "super-lig", "2023-2024", 45270.47916666666, "Istanbulspor AS", 0, "Alanyaspor", 1, 2.61, "03/12/2023 11:42", 3.32, "10/12/2023 11:29", 3.42, "03/12/2023 11:42", 3.7, "10/12/2023 11:26", 2.76, "03/12/2023 11:42", 2.19, "10/12/2023 11:29", "https://www.betexplorer.com/football/turkey/super-lig/istanbulspor-as-alanyaspor/fLtzFd9T/")
    ,(146, "turkey", "super-lig", "2023-2024", 45270.58333333334, "Basaksehir", 1, "Hatayspor", 0, 1.87, "03/12/2023 14:12", 1.91, "10/12/2023 13:58", 3.75, "03/12/2023 14:12", 3.65, "10/12/2023 13:58", 4.17, "03/12/2023 14:12", 4.31, "10/12/2023 13:58", "https://www.betexplorer.com/football/turkey/super-lig/basaksehir-hatayspor/6aAhXwPc/")
    ,(147, "turkey", "super-lig", "2023-2024", 45270.58333333334, "Kayserispor", 2, "Pendikspor", 0, 2.07, "03/12/2023 14:12", 1.69, "10/12/2023 13:57", 3.7, "03/12/2023 14:12", 4.23, "10/12/2023 13:58", 3.5, "03/12/2023 14:12", 4.9, "10/12/2023 13:58", "https://www.betexplorer.com/football/turkey/super-lig/kayserispor-pendikspor/tIxWFxfN/")
    ,(148, "turkey", "super-lig", "2023-2024", 45270.70833333334, "Gaziantep", 1, "Trabzonspor", 3, 2.93, "03/12/2023 17:13", 3.92, "10/12/2023 16:57", 3.32, "03/12/2023 17:13", 3.51, "10/12/2023 16:55", 2.52, "03/12/2023 17:13", 2.04, "10/12/2023 16:55", "https://www.betexplorer.com/football/turkey/super-lig/gaziantep-trabzonspor/S8BlYJ9i/")
    ,(149, "turkey", "super-lig", "2023-2024", 45271.75, "Ankaragucu", 1, "Rizespor", 1, 2.14, "04/12/2023 18:12", 2.04, "11/12/2023 17:58", 3.41, "04/12/2023 18:12", 3.62, "11/12/2023 17:58", 3.6, "04/12/2023 18:12", 3.8, "11/12/2023 17:58", "https://www.betexplorer.com/football/turkey/super-lig/ankaragucu-rizespor/tWH4Uy9G/")
    ,(150, "turkey", "super-lig", "2023-2024", 45271.75, "Konyaspor", 0, "Sivasspor", 1, 2.17, "04/12/2023 18:12", 2.09, "11/12/2023 17:33", 3.4, "04/12/2023 18:12", 3.41, "11/12/2023 17:33", 3.51, "04/12/2023 18:12", 3.9, "11/12/2023 17:33", "https://www.betexplorer.com/football/turkey/super-lig/konyaspor-sivasspor/IeuYC9wq/")
    ,(151, "turkey", "super-lig", "2023-2024", 45279.625, "Basaksehir", 3, "Sivasspor", 1, 2.12, "14/12/2023 09:42", 2.22, "19/12/2023 14:32", 3.35, "14/12/2023 09:42", 3.48, "19/12/2023 14:29", 3.71, "14/12/2023 09:42", 3.47, "19/12/2023 14:57", "https://www.betexplorer.com/football/turkey/super-lig/basaksehir-sivasspor/4WPGTGuK/")
)

$nRows = $newRows.Count
$nCols = 22
$arr = New-Object "object[,]" $nRows,$nCols
for ($r = 0; $r -lt $nRows; $r++) {
    for ($c = 0; $c -lt $nCols; $c++) {
        $arr[$r,$c] = $newRows[$r][$c]
    }
}
$ws.Range("A138:V152").Value2 = $arr

